$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2025-03-17"
$ws.Range("B7").Value = "qasim"
$ws.Range("C7").Value = "Qasim"
$ws.Range("D7").Value = "19:05:54"
$ws.Range("E7").Value = "19:06:16"
$ws.Range("F7").Value = "trainer"
